$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        throw "Could not find text: $old"
    }
    $r = $d.Range($idx, $idx + $old.Length)
    $r.Text = $new
}

# 1. Remove the existing _GoBack bookmark (it will be re-added at its new
#    location once the surrounding text has been updated).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Apply the wording changes from the diff.
Replace-Text "increasingly more distributed" "increasingly distributed"
Replace-Text "hardware nodes" "computational nodes"
Replace-Text "system can provide the application's required services" "system has enough resources to provide the application's required resources"
Replace-Text "integrated our modeling semantics into a" "implemented our modeling semantics in a"

# 3. Re-insert the _GoBack bookmark at its new location (between "in " and "a ").
$full = $d.Content.Text
$idx = $full.IndexOf("implemented our modeling semantics in a")
$insertPos = $idx + ("implemented our modeling semantics in ").Length
$r = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $r)
